# getting familiar with tree traversal
# Adds a new LeetCode entry (day 32 / row 34): "Binary Tree Postorder Traversal"

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# --- New data row (row 34) --------------------------------------------------
$ws1.Range("B34").Value = "Binary Tree Postorder Traversal"
$ws1.Range("C34").Value = 1
$ws1.Range("D34").Value = 1
$ws1.Range("E34").Value = 34
$ws1.Range("F34").Value = 0.19
$ws1.Range("G34").Value = 16.09
$ws1.Range("H34").Value = 0.0217
$ws1.Range("I34").Value = "https://leetcode.com/problems/binary-tree-postorder-traversal/submissions/"

# Hyperlink for the new problem name cell (matches style/pattern of B22:B33).
# NOTE: Hyperlinks.Add's "TextToDisplay" parameter overwrites the cell's
# value, so the problem-name text is re-applied (and the style re-applied,
# since Hyperlinks.Add also resets cell formatting) right after.
$ws1.Hyperlinks.Add($ws1.Range("B34"), "https://leetcode.com/problems/binary-tree-postorder-traversal/", $null, $null, "https://leetcode.com/problems/binary-tree-postorder-traversal/") | Out-Null
$ws1.Range("B34").Value = "Binary Tree Postorder Traversal"
$ws1.Range("B34").Style = "Hyperlink"

# --- Update the selection / active cell on Sheet1 --------------------------
$ws1.Activate()
$ws1.Range("I38").Select() | Out-Null

# --- Sheet2 summary formulas auto-recalculate once data changes ------------
$excel.Calculate()
